$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "198×7=1386" "461×2=922"
Replace-Text "360×2=720" "118×5=590"
Replace-Text "201×4=804" "989×3=2967"
Replace-Text "892×5=4460" "919×4=3676"
Replace-Text "304×5=1520" "951×7=6657"
Replace-Text "392×5=1960" "728×3=2184"
Replace-Text "104×7=728" "163×3=489"
Replace-Text "269×2=538" "512×9=4608"
Replace-Text "852×2=1704" "483×3=1449"
Replace-Text "148×4=592" "762×2=1524"
Replace-Text "507×7=3549" "731×8=5848"
Replace-Text "107×4=428" "748×7=5236"
Replace-Text "819×8=6552" "427×2=854"
Replace-Text "926×2=1852" "224×4=896"
Replace-Text "953×4=3812" "436×6=2616"
Replace-Text "126×4=504" "988×8=7904"
Replace-Text "120×2=240" "609×3=1827"
Replace-Text "245×5=1225" "138×8=1104"
Replace-Text "856×9=7704" "334×7=2338"
Replace-Text "550×4=2200" "622×7=4354"
Replace-Text "516×8=4128" "334×2=668"
Replace-Text "901×2=1802" "791×5=3955"
Replace-Text "672×2=1344" "993×8=7944"
Replace-Text "280×5=1400" "895×5=4475"
Replace-Text "389×9=3501" "856×5=4280"
